function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# 1) Build the new "2022-Q1" sheet by duplicating "2021-Q4" (same header/style layout)
#    and dropping it right before "总计".
$q4.Copy($total)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"
$new.Range("A8:A31").EntireRow.Delete()

# 2) Overwrite the fund-holding detail rows with the 2022-Q1 data
$rows = @(
    @("012930", "中庚价值先锋股票",         "54.59", "94.46", "5.38", "2.9369", 5),
    @("000986", "太平灵活配置混合型发起式", "18.13", "82.45", "4.25", "0.7705", 6),
    @("013004", "国泰价值领航股票A",        "7.51",  "94.17", "3.90", "0.2929", 1),
    @("009537", "太平行业优选股票A",        "0.88",  "90.50", "5.49", "0.0483", 5),
    @("013005", "国泰价值领航股票C",        "0.33",  "94.17", "3.90", "0.0129", 1),
    @("009538", "太平行业优选股票C",        "0.20",  "90.50", "5.49", "0.0110", 5)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2
    Set-TextValue $new.Cells.Item($rowNum, 2) $r[0]
    $new.Cells.Item($rowNum, 3).Value = $r[1]
    Set-TextValue $new.Cells.Item($rowNum, 4) $r[2]
    Set-TextValue $new.Cells.Item($rowNum, 5) $r[3]
    Set-TextValue $new.Cells.Item($rowNum, 6) $r[4]
    Set-TextValue $new.Cells.Item($rowNum, 7) $r[5]
    $new.Cells.Item($rowNum, 8).Value = $r[6]
}
